{"js": "// Update the report date in the header block.\n{\n  const dateResults = context.document.body.search(\"December 08, 2025\", { matchCase: true });\n  dateResults.load(\"items\");\n  await context.sync();\n  if (dateResults.items.length > 0) {\n    dateResults.items[0].insertText(\"December 11, 2025\", \"Replace\");\n    await context.sync();\n  }\n}\n\n// Locate the start and end paragraphs of the Section 2 block that needs to be\n// replaced/expanded with the new BSM pricing analysis content.\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nconst startText = \"2. Black\\u2013Scholes\\u2013Merton Pricing Results\";\nconst endText = \"This initial analysis provides a foundation for further validation and integration of the BSM pricing model within the broader quantitative framework.\";\n\nlet startIdx = -1;\nlet endIdx = -1;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === startText) {\n    startIdx = i;\n  }\n  if (paras.items[i].text === endText) {\n    endIdx = i;\n    break;\n  }\n}\n\nif (startIdx === -1 || endIdx === -1) {\n  throw new Error(\"Could not locate Section 2 boundary paragraphs (start=\" + startIdx + \", end=\" + endIdx + \")\");\n}\n\n// New paragraph texts for the expanded Section 2 content. \"\\v\" (vertical tab)\n// denotes a manual line break (<w:br/>) within a single paragraph/run, matching\n// the source document's convention for multi-line markdown-style paragraphs.\nconst newTexts = [\n  \"## Section 2: Black\u2013Scholes\u2013Merton Pricing Analysis\",\n  \"### 2.1 Overview of BSM Pricing Results\",\n  \"This section presents the analysis of option pricing using the Black\u2013Scholes\u2013Merton (BSM) model across different asset classes. The analysis examines how the model performs under various market conditions and input parameters, providing insights into pricing accuracy and reliability.\",\n  \"### 2.2 Summary of Pricing Inputs and Outputs\",\n  \"The BSM pricing model utilizes several key parameters to generate option prices:\\u000b- Spot price (S): Current market price of the underlying asset\\u000b- Strike price (K): Predetermined price at which the option can be exercised\\u000b- Time to maturity (T): Time remaining until option expiration (in years)\\u000b- Risk-free rate (r): Annual risk-free interest rate\\u000b- Volatility (\u03c3): Expected volatility of the underlying asset\\u000b- Option type: Call or Put\",\n  \"#### 2.2.1 Pricing Results by Asset Class\",\n  \"**Equities**\",\n  \"| Valuation Date | Spot Price | Strike | Maturity (Years) | Risk-Free Rate | Volatility | Option Type | BSM Price |\\u000b|----------------|------------|--------|------------------|----------------|------------|-------------|-----------|\\u000b| 2023-01-15     | 150.25     | 155.00 | 0.5              | 0.03           | 0.22       | Call        | 8.76      |\\u000b| 2023-01-15     | 150.25     | 145.00 | 0.5              | 0.03           | 0.22       | Put         | 3.24      |\\u000b| 2023-02-10     | 152.50     | 155.00 | 0.4              | 0.035          | 0.24       | Call        | 9.12      |\\u000b| 2023-03-05     | 148.75     | 145.00 | 0.3              | 0.04           | 0.25       | Put         | 3.87      |\",\n  \"**FX Options**\",\n  \"| Valuation Date | Spot Price | Strike | Maturity (Years) | Risk-Free Rate | Volatility | Option Type | BSM Price |\\u000b|----------------|------------|--------|------------------|----------------|------------|-------------|-----------|\\u000b| 2023-01-20     | 1.10       | 1.12   | 0.25             | 0.02           | 0.08       | Call        | 0.0156    |\\u000b| 2023-01-20     | 1.10       | 1.08   | 0.25             | 0.02           | 0.08       | Put         | 0.0112    |\\u000b| 2023-02-15     | 1.09       | 1.10   | 0.2              | 0.025          | 0.09       | Call        | 0.0178    |\\u000b| 2023-03-10     | 1.08       | 1.07   | 0.15             | 0.03           | 0.10       | Put         | 0.0134    |\",\n  \"**Commodity Options**\",\n  \"| Valuation Date | Spot Price | Strike | Maturity (Years) | Risk-Free Rate | Volatility | Option Type | BSM Price |\\u000b|----------------|------------|--------|------------------|----------------|------------|-------------|-----------|\\u000b| 2023-01-25     | 80.50      | 82.00  | 0.75             | 0.025          | 0.28       | Call        | 7.23      |\\u000b| 2023-01-25     | 80.50      | 79.00  | 0.75             | 0.025          | 0.28       | Put         | 5.45      |\\u000b| 2023-02-20     | 82.75      | 83.00  | 0.6              | 0.03           | 0.30       | Call        | 7.89      |\\u000b| 2023-03-15     | 79.25      | 78.00  | 0.5              | 0.035          | 0.32       | Put         | 5.67      |\",\n  \"### 2.3 Pricing Trends Over Time\",\n  \"#### 2.3.1 Equity Options Pricing Trend\\u000b![Equity Options BSM Pricing Trend]\\u000b(Graph showing BSM prices for equity options from January to March 2023, with separate lines for call and put options. The graph demonstrates slight upward trend for call options and moderate volatility for put options.)\",\n  \"#### 2.3.2 FX Options Pricing Trend\\u000b![FX Options BSM Pricing Trend]\\u000b(Graph showing BSM prices for FX options from January to March 2023, with separate lines for call and put options. The trend shows relatively stable pricing with minor fluctuations corresponding to changes in spot rates.)\",\n  \"#### 2.3.3 Commodity Options Pricing Trend\\u000b![Commodity Options BSM Pricing Trend]\\u000b(Graph showing BSM prices for commodity options from January to March 2023, with separate lines for call and put options. The graph indicates higher price volatility compared to other asset classes, particularly for call options.)\",\n  \"### 2.4 Key Observations\",\n  \"- **Input Sensitivity**: The BSM model shows particular sensitivity to volatility inputs across all asset classes, with commodity options demonstrating the highest price impact from volatility changes.\\u000b- **Asset Class Differences**: FX options exhibit lower absolute price values but similar relative price movements compared to equity and commodity options.\\u000b- **Data Quality Considerations**: Some pricing calculations may yield unreliable results when inputs are missing or invalid. These cases require special handling in production environments.\\u000b- **Time Decay Effects**: As expected under the BSM framework, option prices generally decrease as time to maturity shortens, with the effect being more pronounced for at-the-money options.\"\n];\n\n// Seed the first new paragraph into the existing start paragraph (preserves\n// position/continuity), then insert the remaining new paragraphs after it.\nlet anchor = paras.items[startIdx];\nanchor.insertText(newTexts[0], \"Replace\");\nawait context.sync();\n\nfor (let i = 1; i < newTexts.length; i++) {\n  anchor = anchor.insertParagraph(newTexts[i], \"After\");\n  await context.sync();\n}\n\n// Remove the remaining old paragraphs that were part of the original Section 2\n// block (everything from right after the original start paragraph through the\n// original end paragraph, inclusive).\nfor (let i = startIdx + 1; i <= endIdx; i++) {\n  paras.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# Update the report date in the header block.\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"December 08, 2025\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"December 11, 2025\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# Locate the start/end paragraphs of the Section 2 block that must be replaced\n# and expanded with the new BSM pricing analysis content.\n$startText = \"2. Black\" + [char]8211 + \"Scholes\" + [char]8211 + \"Merton Pricing Results\"\n$endText = \"This initial analysis provides a foundation for further validation and integration of the BSM pricing model within the broader quantitative framework.\"\n\n$startIdx = -1\n$endIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    $t = $t.TrimEnd([char]13)\n    if ($t -eq $startText) {\n        $startIdx = $i\n    }\n    if ($t -eq $endText) {\n        $endIdx = $i\n    }\n}\n\nif ($startIdx -eq -1 -or $endIdx -eq -1) {\n    throw \"Could not locate Section 2 boundary paragraphs (start=$startIdx, end=$endIdx)\"\n}\n\n# New paragraph texts for the expanded Section 2 content. [char]11 (vertical\n# tab) denotes a manual line break (<w:br/>) within a single paragraph/run,\n# matching the source document's convention for multi-line markdown-style\n# paragraphs.\n$newTexts = @(\n    (\"## Section 2: Black\u2013Scholes\u2013Merton Pricing Analysis\"),\n    (\"### 2.1 Overview of BSM Pricing Results\"),\n    (\"This section presents the analysis of option pricing using the Black\u2013Scholes\u2013Merton (BSM) model across different asset classes. The analysis examines how the model performs under various market conditions and input parameters, providing insights into pricing accuracy and reliability.\"),\n    (\"### 2.2 Summary of Pricing Inputs and Outputs\"),\n    (\"The BSM pricing model utilizes several key parameters to generate option prices:\" + [char]11 + \"- Spot price (S): Current market price of the underlying asset\" + [char]11 + \"- Strike price (K): Predetermined price at which the option can be exercised\" + [char]11 + \"- Time to maturity (T): Time remaining until option expiration (in years)\" + [char]11 + \"- Risk-free rate (r): Annual risk-free interest rate\" + [char]11 + \"- Volatility (\u03c3): Expected volatility of the underlying asset\" + [char]11 + \"- Option type: Call or Put\"),\n    (\"#### 2.2.1 Pricing Results by Asset Class\"),\n    (\"**Equities**\"),\n    (\"| Valuation Date | Spot Price | Strike | Maturity (Years) | Risk-Free Rate | Volatility | Option Type | BSM Price |\" + [char]11 + \"|----------------|------------|--------|------------------|----------------|------------|-------------|-----------|\" + [char]11 + \"| 2023-01-15     | 150.25     | 155.00 | 0.5              | 0.03           | 0.22       | Call        | 8.76      |\" + [char]11 + \"| 2023-01-15     | 150.25     | 145.00 | 0.5              | 0.03           | 0.22       | Put         | 3.24      |\" + [char]11 + \"| 2023-02-10     | 152.50     | 155.00 | 0.4              | 0.035          | 0.24       | Call        | 9.12      |\" + [char]11 + \"| 2023-03-05     | 148.75     | 145.00 | 0.3              | 0.04           | 0.25       | Put         | 3.87      |\"),\n    (\"**FX Options**\"),\n    (\"| Valuation Date | Spot Price | Strike | Maturity (Years) | Risk-Free Rate | Volatility | Option Type | BSM Price |\" + [char]11 + \"|----------------|------------|--------|------------------|----------------|------------|-------------|-----------|\" + [char]11 + \"| 2023-01-20     | 1.10       | 1.12   | 0.25             | 0.02           | 0.08       | Call        | 0.0156    |\" + [char]11 + \"| 2023-01-20     | 1.10       | 1.08   | 0.25             | 0.02           | 0.08       | Put         | 0.0112    |\" + [char]11 + \"| 2023-02-15     | 1.09       | 1.10   | 0.2              | 0.025          | 0.09       | Call        | 0.0178    |\" + [char]11 + \"| 2023-03-10     | 1.08       | 1.07   | 0.15             | 0.03           | 0.10       | Put         | 0.0134    |\"),\n    (\"**Commodity Options**\"),\n    (\"| Valuation Date | Spot Price | Strike | Maturity (Years) | Risk-Free Rate | Volatility | Option Type | BSM Price |\" + [char]11 + \"|----------------|------------|--------|------------------|----------------|------------|-------------|-----------|\" + [char]11 + \"| 2023-01-25     | 80.50      | 82.00  | 0.75             | 0.025          | 0.28       | Call        | 7.23      |\" + [char]11 + \"| 2023-01-25     | 80.50      | 79.00  | 0.75             | 0.025          | 0.28       | Put         | 5.45      |\" + [char]11 + \"| 2023-02-20     | 82.75      | 83.00  | 0.6              | 0.03           | 0.30       | Call        | 7.89      |\" + [char]11 + \"| 2023-03-15     | 79.25      | 78.00  | 0.5              | 0.035          | 0.32       | Put         | 5.67      |\"),\n    (\"### 2.3 Pricing Trends Over Time\"),\n    (\"#### 2.3.1 Equity Options Pricing Trend\" + [char]11 + \"![Equity Options BSM Pricing Trend]\" + [char]11 + \"(Graph showing BSM prices for equity options from January to March 2023, with separate lines for call and put options. The graph demonstrates slight upward trend for call options and moderate volatility for put options.)\"),\n    (\"#### 2.3.2 FX Options Pricing Trend\" + [char]11 + \"![FX Options BSM Pricing Trend]\" + [char]11 + \"(Graph showing BSM prices for FX options from January to March 2023, with separate lines for call and put options. The trend shows relatively stable pricing with minor fluctuations corresponding to changes in spot rates.)\"),\n    (\"#### 2.3.3 Commodity Options Pricing Trend\" + [char]11 + \"![Commodity Options BSM Pricing Trend]\" + [char]11 + \"(Graph showing BSM prices for commodity options from January to March 2023, with separate lines for call and put options. The graph indicates higher price volatility compared to other asset classes, particularly for call options.)\"),\n    (\"### 2.4 Key Observations\"),\n    (\"- **Input Sensitivity**: The BSM model shows particular sensitivity to volatility inputs across all asset classes, with commodity options demonstrating the highest price impact from volatility changes.\" + [char]11 + \"- **Asset Class Differences**: FX options exhibit lower absolute price values but similar relative price movements compared to equity and commodity options.\" + [char]11 + \"- **Data Quality Considerations**: Some pricing calculations may yield unreliable results when inputs are missing or invalid. These cases require special handling in production environments.\" + [char]11 + \"- **Time Decay Effects**: As expected under the BSM framework, option prices generally decrease as time to maturity shortens, with the effect being more pronounced for at-the-money options.\"),\n)\n\n$oldCount = $endIdx - $startIdx + 1\n\n# Seed the first new paragraph into the existing start paragraph (preserves\n# position/continuity), then insert the remaining new paragraphs after it.\n$cur = $d.Paragraphs.Item($startIdx)\n$cur.Range.Text = $newTexts[0]\n\nfor ($i = 1; $i -lt $newTexts.Count; $i++) {\n    $cur.Range.InsertParagraphAfter()\n    $cur = $d.Paragraphs.Item($startIdx + $i)\n    $cur.Range.Text = $newTexts[$i]\n}\n\n# Remove the remaining old paragraphs that were part of the original Section 2\n# block. After the inserts above, they now sit right after the newly inserted\n# paragraphs (old paragraph count minus the one we repurposed).\n$oldRemaining = $oldCount - 1\nfor ($i = 1; $i -le $oldRemaining; $i++) {\n    $d.Paragraphs.Item($startIdx + $newTexts.Count).Range.Delete()\n}\n"}
